$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# ---------------------------------------------------------------------------
# Hoja1
# ---------------------------------------------------------------------------

# Header "Importe" becomes a text column (numeric amounts stored as text)
$ws1.Range("B1").NumberFormat = "@"
$ws1.Range("B2:B15").NumberFormat = "@"

# Column A (# Empleado) widened
$ws1.Columns.Item(1).ColumnWidth = 25.59
$ws1.Columns.Item(1).EntireColumn.AutoFit() | Out-Null

# Employee numbers (column A)
$aVals = @(280,281,282,283,299,286,287,284,290,292,293,294,295,304)
for ($i = 0; $i -lt $aVals.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 1).Value = $aVals[$i]
}

# Amounts (column B, stored as text) and periodicity (column C)
$bVals = @("185000.10","70684","115414","62914.20","24226.20","31700.20","34500.10","322800.10","50500.20","23140.79","24226.20","42222.49","24226.20","24000")
for ($i = 0; $i -lt $bVals.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $bVals[$i]
    $ws1.Cells.Item($row, 3).Value = "Mensual"
}

# Rows 7-15 are brand new cells in column C; clear the style picked up from
# neighbouring column D so they stay unformatted like C2:C6.
$ws1.Range("C7:C15").Style = "Normal"

# Active cell moved to F7 on this sheet
$ws1.Range("F7").Select()

# ---------------------------------------------------------------------------
# Hoja2 - no value change, cosmetic/re-saved formatting only
# ---------------------------------------------------------------------------
